{"js": "// Replace the date line and each multiplication-problem cell with its new value.\n// Every old string below occurs exactly once in the document, so an exact-text\n// search (matchCase, no wildcards) safely targets the single run to update.\nconst replacements = [\n  [\"2024-01-29 Monday\", \"2024-01-30 Tuesday\"],\n  [\"153\u00d73=\", \"618\u00d78=\"],\n  [\"885\u00d77=\", \"588\u00d74=\"],\n  [\"164\u00d76=\", \"917\u00d79=\"],\n  [\"547\u00d79=\", \"884\u00d78=\"],\n  [\"917\u00d75=\", \"799\u00d76=\"],\n  [\"876\u00d79=\", \"601\u00d78=\"],\n  [\"462\u00d77=\", \"855\u00d76=\"],\n  [\"641\u00d78=\", \"827\u00d77=\"],\n  [\"579\u00d75=\", \"354\u00d75=\"],\n  [\"343\u00d74=\", \"886\u00d78=\"],\n  [\"591\u00d77=\", \"559\u00d77=\"],\n  [\"808\u00d79=\", \"451\u00d74=\"],\n  [\"537\u00d78=\", \"647\u00d75=\"],\n  [\"798\u00d72=\", \"330\u00d77=\"],\n  [\"996\u00d77=\", \"983\u00d73=\"],\n  [\"896\u00d79=\", \"328\u00d77=\"],\n  [\"550\u00d74=\", \"674\u00d78=\"],\n  [\"369\u00d76=\", \"456\u00d79=\"],\n  [\"424\u00d73=\", \"145\u00d74=\"],\n  [\"605\u00d77=\", \"604\u00d73=\"],\n  [\"509\u00d79=\", \"188\u00d72=\"],\n  [\"544\u00d78=\", \"773\u00d79=\"],\n  [\"321\u00d79=\", \"218\u00d73=\"],\n  [\"582\u00d73=\", \"881\u00d73=\"],\n  [\"208\u00d75=\", \"564\u00d72=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Text not found: ' + oldText);\n  }\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();", "ps1": "# Word COM (PowerShell-style) script: update the worksheet date line and\n# every three-digit x one-digit multiplication prompt to the new values.\n# Each \"old\" string is unique in the document, so Find/Replace (wdReplaceAll,\n# scoped by the exact FindText) safely retargets exactly one run per call.\n$d = $word.ActiveDocument\n\n$olds = @(\"2024-01-29 Monday\", \"153\u00d73=\", \"885\u00d77=\", \"164\u00d76=\", \"547\u00d79=\", \"917\u00d75=\", \"876\u00d79=\", \"462\u00d77=\", \"641\u00d78=\", \"579\u00d75=\", \"343\u00d74=\", \"591\u00d77=\", \"808\u00d79=\", \"537\u00d78=\", \"798\u00d72=\", \"996\u00d77=\", \"896\u00d79=\", \"550\u00d74=\", \"369\u00d76=\", \"424\u00d73=\", \"605\u00d77=\", \"509\u00d79=\", \"544\u00d78=\", \"321\u00d79=\", \"582\u00d73=\", \"208\u00d75=\")\n$news = @(\"2024-01-30 Tuesday\", \"618\u00d78=\", \"588\u00d74=\", \"917\u00d79=\", \"884\u00d78=\", \"799\u00d76=\", \"601\u00d78=\", \"855\u00d76=\", \"827\u00d77=\", \"354\u00d75=\", \"886\u00d78=\", \"559\u00d77=\", \"451\u00d74=\", \"647\u00d75=\", \"330\u00d77=\", \"983\u00d73=\", \"328\u00d77=\", \"674\u00d78=\", \"456\u00d79=\", \"145\u00d74=\", \"604\u00d73=\", \"188\u00d72=\", \"773\u00d79=\", \"218\u00d73=\", \"881\u00d73=\", \"564\u00d72=\")\n\nfor ($i = 0; $i -lt $olds.Count; $i++) {\n    $old = $olds[$i]\n    $new = $news[$i]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
